$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values per diff
$ws.Range("F2").Value = -225
$ws.Range("F3").Value = -422
$ws.Range("F4").Value = -709
$ws.Range("B5").Value = "DC_001"
$ws.Range("F5").Value = -98
$ws.Range("H5").Value = 4

# Add new row 6: MAT_B, DC_002, 45293, Distribution Demand - Forecast, 1, -48, 45292, 1
$ws.Range("A6").Value = "MAT_B"
$ws.Range("B6").Value = "DC_002"
$ws.Range("C6").Value = 45293
$ws.Range("C6").NumberFormat = $ws.Range("C5").NumberFormat
$ws.Range("D6").Value = "Distribution Demand - Forecast"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = -48
$ws.Range("G6").Value = 45292
$ws.Range("G6").NumberFormat = $ws.Range("G5").NumberFormat
$ws.Range("H6").Value = 1

# Add new row 7: MAT_B, PLANT_001, 45293, Distribution Demand - Forecast, 0, -98, 45292, 1
$ws.Range("A7").Value = "MAT_B"
$ws.Range("B7").Value = "PLANT_001"
$ws.Range("C7").Value = 45293
$ws.Range("C7").NumberFormat = $ws.Range("C5").NumberFormat
$ws.Range("D7").Value = "Distribution Demand - Forecast"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = -98
$ws.Range("G7").Value = 45292
$ws.Range("G7").NumberFormat = $ws.Range("G5").NumberFormat
$ws.Range("H7").Value = 1
